# Insert a new "Any Questions?" title slide right before the final
# "Acknowledgments" slide (i.e. at position 11 in the now-12-slide deck).
$p = $ppt.ActivePresentation

$count = $p.Slides.Count
$newSlide = $p.Slides.Add($count, 1)   # 1 = ppLayoutTitle ("Title Slide" -> ctrTitle/subTitle)

# Title placeholder: "          Any Questions?" followed by a blank line.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "          Any Questions?`r"

# Subtitle placeholder is left empty (matches the authored slide).
